# [Doc] Add problem statement slide to ppt
#
# 1. Insert a new "Problem Statement" slide right after the title slide.
# 2. Tweak wording on the Exploratory Data Analysis (EDA) slide.
# 3. Refresh the cached "today" date field shown on every slide footer
#    (slide master + every slide layout) from 6/1/2025 to 17/1/2025.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper: find a slide's 1-based index by its title text.
# ---------------------------------------------------------------------
function Get-SlideIndexByTitle($title) {
    for ($i = 1; $i -le $p.Slides.Count; $i++) {
        $s = $p.Slides.Item($i)
        if ($s.Shapes.HasTitle) {
            if ($s.Shapes.Title.TextFrame.TextRange.Text -eq $title) {
                return $i
            }
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) Insert the new "Problem Statement" slide at position 2 (right
#    after the title slide), using the "Title and Content" layout.
# ---------------------------------------------------------------------
$titleAndContent = 2
$newSlide = $p.Slides.Add(2, $titleAndContent)

$newTitle = $newSlide.Shapes.Item(1)
$newTitle.TextFrame.TextRange.Text = "Problem Statement"

$newBody = $newSlide.Shapes.Item(2)
$bodyTr = $newBody.TextFrame.TextRange

$line1 = "AgroTech innovations faces challenges in optimizing crop yields and resource management due to inefficiencies"
$line2 = "Goal:"
$line3 = "Predict temperature conditions within farms" + [char]0x2019 + " closed environment to ensure optimal plant growth"
$line4 = "Categorize combined " + [char]0x201C + "Plant Type-Stage" + [char]0x201D + " based on sensor data to aid in strategic planning and resource allocation"

$bodyTr.Text = $line1 + "`r" + $line2 + "`r" + $line3 + "`r" + $line4

# NOTE: reading TextRange.Text back normalizes "smart" punctuation
# (curly quotes/apostrophes) to their ASCII look-alikes, even though the
# underlying/saved XML keeps the real Unicode characters. So offsets for
# the later paragraphs must be computed from the known source lengths
# instead of searching inside the read-back text.
$idx3 = $line1.Length + 1 + $line2.Length + 1
$idx4 = $idx3 + $line3.Length + 1

$sub3 = $bodyTr.Characters($idx3 + 1, $line3.Length)
$sub3.IndentLevel = 2
$sub4 = $bodyTr.Characters($idx4 + 1, $line4.Length)
$sub4.IndentLevel = 2

# ---------------------------------------------------------------------
# 2) Wording tweaks on the "Exploratory Data Analysis (EDA)" slide.
# ---------------------------------------------------------------------
$edaIdx = Get-SlideIndexByTitle "Exploratory Data Analysis (EDA)"
if ($edaIdx -gt 0) {
    $edaSlide = $p.Slides.Item($edaIdx)
    for ($j = 1; $j -le $edaSlide.Shapes.Count; $j++) {
        $shp = $edaSlide.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -like "*Part-specific*") {
                    $full = $tr.Text
                    $full = $full.Replace("Part-specific", "Task-specific (Label-encoded)")
                    $oldLine = "Remaining columns with " + [char]0x2018 + "object" + [char]0x2019 + " type data are one-hot encoded"
                    $newLine = "Columns with " + [char]0x2018 + "object" + [char]0x2019 + " type data are one-hot encoded"
                    $full = $full.Replace($oldLine, $newLine)
                    $tr.Text = $full
                }
            }
        }
    }
}

# ---------------------------------------------------------------------
# 3) Refresh cached date field text (6/1/2025 -> 17/1/2025) on the
#    slide master and every slide layout's Date Placeholder.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $t = $shp.TextFrame.TextRange.Text
                if ($t -eq "6/1/2025") {
                    $shp.TextFrame.TextRange.Text = "17/1/2025"
                }
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}
